# Move the NN-specific columns (layers / metrics / optimizer / fit) off the
# "Classifier" sheet and onto the "ClassifierParams" sheet, carrying their
# data validation along with them.

$wb = $excel.ActiveWorkbook

$classifier = $wb.Worksheets.Item("Classifier")
$classifierParams = $wb.Worksheets.Item("ClassifierParams")

# --- Remove H1:K1 from Classifier (clears header text + any data validation
# that was anchored to those columns, e.g. the optimizer list on J2:J1048576)
$classifier.Range("H1:K1").Clear()
$classifier.Columns("J").Validation.Delete()

# --- Add the corresponding headers to ClassifierParams (C1:F1)
$classifierParams.Range("C1").Value = "classifierParams__layers"
$classifierParams.Range("D1").Value = "classifierParams__metrics"
$classifierParams.Range("E1").Value = "classifierParams__optimizer"
$classifierParams.Range("F1").Value = "classifierParams__fit"

# --- Add the optimizer list validation to column E (E2:E1048576) on
# ClassifierParams, matching what used to live on Classifier!J2:J1048576
$optRange = $classifierParams.Range("E2:E1048576")
$optRange.Validation.Delete()
$optRange.Validation.Add(3, 1, 1, '"Adagrad,Adam,Adamax,Nadam,SGD"')
$optRange.Validation.IgnoreBlank = $true
$optRange.Validation.InCellDropdown = $true
$optRange.Validation.ShowInput = $true
$optRange.Validation.ShowError = $true
